$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 271
$ws.Range("C3").Value = 168607
$ws.Range("C4").Value = 159460
$ws.Range("C5").Value = 9147
$ws.Range("C7").Value = 5.43
$ws.Range("C8").Value = 65.38
